# Settings.xlsx update: "Added all reports forms to the repository."
#
# The four "technical" parameter rows (ReportsPath, RawDataPath, Parameters,
# ReportsPrepared) get their human readable B-column values replaced with
# numbered Russian folder names. Row order in which the new values are
# written matters because it controls the order new entries are appended
# to the shared-strings table (B6, B2, B4, B3 - matching the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "01. Сформированные отчеты"
$ws.Range("B2").Value = "03. Формы отчетов"
$ws.Range("B4").Value = "04. Настройки"
$ws.Range("B3").Value = "02. Данные выгруженные из DES.LM"

# Move/restore the active selection to B3.
$ws.Range("B3").Select()

# Page setup for printing: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Strip personal/author information on save (workbookPr filterPrivacy).
$wb.RemovePersonalInformation = $true
